# Updated cryptos list on Sat May 25 02:50:43 UTC 2024 with GitHub Actions
#
# Refreshes the "cryptos" worksheet with the latest scraped Price /
# Volume(1h) snapshot for each coin row (A2:E51). Also reflects that
# VeChain and Maker swapped rank positions (rows 50 and 51) in this run.
#
# Numeric-looking price strings (e.g. "600.41") are written with a
# temporary Text number format so Excel keeps them as literal strings
# (matching the source data, which uses "."-grouped numbers) instead of
# auto-converting them to the Number type; the cell format is restored
# to Normal/General immediately afterwards so no visible formatting
# changes are introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.737.25'
$ws.Range("E2").Value = '  +1.50%  '
$ws.Range("D3").Value = '3.727.05'
$ws.Range("E3").Value = '  -2.02%  '
$ws.Range("E4").Value = '  -0.01%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '600.41'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E5").Value = '  -0.19%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '169.11'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E6").Value = '  -1.73%  '
$ws.Range("D7").Value = '3.725.70'
$ws.Range("E7").Value = '  -2.05%  '
$ws.Range("E8").Value = '  +0.04%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.535'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E9").Value = '  +0.59%  '
$ws.Range("E10").Value = '  +2.56%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '6.33'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E11").Value = '  +1.88%  '
$ws.Range("E12").Value = '  -1.06%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '38.12'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E13").Value = '  -1.41%  '
$ws.Range("E14").Value = '  +0.61%  '
$ws.Range("D15").Value = '4.347.09'
$ws.Range("E15").Value = '  -1.87%  '
$ws.Range("D16").Value = '3.727.32'
$ws.Range("E16").Value = '  -1.87%  '
$ws.Range("D17").Value = '68.736.91'
$ws.Range("E17").Value = '  +1.57%  '
$ws.Range("E18").Value = '  +0.58%  '
$ws.Range("E19").Value = '  +1.00%  '
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '17.15'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E20").Value = '  -0.05%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '493.45'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E21").Value = '  +0.41%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '10.59'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E22").Value = '  +15.66%  '
$ws.Range("E23").Value = '  -2.24%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '85.04'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E24").Value = '  -1.00%  '
$ws.Range("E25").Value = '  -0.92%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '2.31'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E26").Value = '  -2.95%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '12.44'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E27").Value = '  +0.85%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '10.14'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E28").Value = '  -1.01%  '
$ws.Range("E29").Value = '  +0.01%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '2.57'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E30").Value = '  +5.76%  '
$ws.Range("E31").Value = '  -0.46%  '
$ws.Range("E32").Value = '  +1.80%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '31.47'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E33").Value = '  -3.93%  '
$ws.Range("D34").Value = '3.867.68'
$ws.Range("E35").Value = '  -0.87%  '
$ws.Range("D36").Value = '3.656.37'
$ws.Range("E36").Value = '  -2.07%  '
$ws.Range("E37").Value = '  +0.01%  '
$ws.Range("E38").Value = '  -0.40%  '
$ws.Range("E39").Value = '  +0.26%  '
$ws.Range("E40").Value = '  -0.35%  '
$ws.Range("E41").Value = '  -1.29%  '
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '437.47'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E42").Value = '  -5.21%  '
$ws.Range("E43").Value = '  -0.51%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '1.98'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E44").Value = '  -1.33%  '
$ws.Range("E45").Value = '  +1.04%  '
$ws.Range("E46").Value = '  +0.79%  '
$ws.Range("E47").Value = '  +0.02%  '
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '40.62'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E48").Value = '  -0.88%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '141.18'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E49").Value = '  +1.65%  '
$ws.Range("B50").Value = 'Maker'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D50").Value = '2.774.26'
$ws.Range("E50").Value = '  -2.60%  '
$ws.Range("B51").Value = 'VeChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.0354'
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E51").Value = '  +0.73%  '
